$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 12, shifting existing rows 12:89 down to 13:90
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new weekly record
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value = "Maule"
$ws.Cells.Item(12, 4).Value = 44901
$ws.Cells.Item(12, 5).Value = 7
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100101
$ws.Cells.Item(12, 8).Value = "Berries"
$ws.Cells.Item(12, 9).Value = 100101001
$ws.Cells.Item(12, 10).Value = "Arándano (blue)"
$ws.Cells.Item(12, 11).Value = "Sin especificar"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 230
$ws.Cells.Item(12, 14).Value = 3000
$ws.Cells.Item(12, 15).Value = 3000
$ws.Cells.Item(12, 16).Value = 3000
$ws.Cells.Item(12, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(12, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(12, 19).Value = 1500
$ws.Cells.Item(12, 20).Value = 2
